$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Insert 3 new rows at the top. All existing data (rows 1-32, plus the 4
# trailing blank rows) shifts down by 3 -- this also carries formulas,
# shared-string references and styles along for free.
# ---------------------------------------------------------------------------
$ws.Rows("1:3").Insert()

# Copy the (now shifted) first data row's formatting onto the 3 new rows so
# they pick up the same cell styles (date format on A, text formats on
# D/F/G) without minting brand new style entries.
$ws.Range("A4:G4").Copy()
$ws.Range("A1:G3").PasteSpecial(-4122)   # xlPasteFormats

# Trailing characters on the monto (F) column values are non-breaking
# spaces (U+00A0 U+00A0), matching the existing shared-string entries for
# the other "monto" values in this sheet (e.g. "0.75\xA0\xA0").
$nbsp2 = [string]([char]0x00A0) + [string]([char]0x00A0)

# ---------------------------------------------------------------------------
# New rows 1-3 data. Shared-string table entries get minted in the order
# cells are written, and the canonical edit's sharedStrings.xml appends the
# three new "documento" (D) values before the three new "saldo" (G) values
# -- so write column-by-column (A,B,C,D,E,F first, then G) to reproduce the
# same new shared-string ordering (70-75).
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = 41684
$ws.Range("B1").Value = "INTERES A SU FAVOR"
$ws.Range("C1").Value = "C"
$ws.Range("D1").Value = "0000950378"
$ws.Range("E1").Value = "AGENCIA PARA PROCESOS BATCH"
$ws.Range("F1").Value = "0.75" + $nbsp2

$ws.Range("A2").Value = 41683
$ws.Range("B2").Value = "INTERES A SU FAVOR"
$ws.Range("C2").Value = "C"
$ws.Range("D2").Value = "0000950389"
$ws.Range("E2").Value = "AGENCIA PARA PROCESOS BATCH"
$ws.Range("F2").Value = "0.25" + $nbsp2

$ws.Range("A3").Value = 41682
$ws.Range("B3").Value = "INTERES A SU FAVOR"
$ws.Range("C3").Value = "C"
$ws.Range("D3").Value = "0000950666"
$ws.Range("E3").Value = "AGENCIA PARA PROCESOS BATCH"
$ws.Range("F3").Value = "0.25" + $nbsp2

$ws.Range("G1").Value = "4015.55"
$ws.Range("G2").Value = "4014.80"
$ws.Range("G3").Value = "4014.55"

# ---------------------------------------------------------------------------
# H column formulas - also fix the small PHP-array typo (missing comma
# before 'mo_borrado_logico') and the date format mask used for
# mo_fecha_crea, same as the canonical edit.
# ---------------------------------------------------------------------------
$formula1 = '=CONCATENATE("array(''mo_fecha'' => new \DateTime(''",TEXT(A1,"yyyy-mm-dd"),"''), ''mo_concepto'' => ''",B1,"'', ''mo_tipo'' => ''",C1,"'', ''mo_documento'' => ''",D1,"'', ''mo_oficina'' => ''",E1,"'', ''mo_monto'' => ",F1,", ''mo_saldo'' => ",G1,", ''mo_fecha_crea'' => new \DateTime(''",TEXT(NOW(),"yyyy-mm-dd HH:mm:ss"),"''), ''mo_quien_crea'' => 1, ''mo_fecha_modifica'' => NULL, ''mo_quien_modifica'' => NULL, ''mo_borrado_logico'' => false),")'
$ws.Range("H1").Formula = $formula1

$formula2 = '=CONCATENATE("array(''mo_fecha'' => new \DateTime(''",TEXT(A2,"yyyy-mm-dd"),"''), ''mo_concepto'' => ''",B2,"'', ''mo_tipo'' => ''",C2,"'', ''mo_documento'' => ''",D2,"'', ''mo_oficina'' => ''",E2,"'', ''mo_monto'' => ",F2,", ''mo_saldo'' => ",G2,", ''mo_fecha_crea'' => new \DateTime(''",TEXT(NOW(),"yyyy-mm-dd HH:mm:ss"),"''), ''mo_quien_crea'' => 1, ''mo_fecha_modifica'' => NULL, ''mo_quien_modifica'' => NULL, ''mo_borrado_logico'' => false),")'
$ws.Range("H2:H3").Formula = $formula2

# ---------------------------------------------------------------------------
# Sheet metadata: dimension grows to A1:H39 (handled automatically by the
# row insert + new cell writes above); update the visible selection.
# ---------------------------------------------------------------------------
$ws.Range("H1:H3").Select()
